{"js": "// Add the new \"to-do\" lines after \"Ajouter doc avec licence et nos noms\",\n// but before the existing trailing empty paragraph.\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\n// Anchor: the existing last paragraph of the body (an empty paragraph).\nconst lastParagraph = paragraphs.getLast();\n\nconst newLines = [\n  \"Mettre une scrollbarre sur la page principale de l\\u2019app\",\n  \"Modifier la livraison quand on modifie le produit\",\n  \"\",\n  \"Questions : \",\n  \"\",\n  \"Creer objet prime\",\n  \"Creer objet arret maladie\",\n  \"Page avec tous les employes et on peut completer avec les salaires qui ont ete modifies\",\n  \"Charges ou les impactees\",\n  \"\"\n];\n\nfor (const line of newLines) {\n  lastParagraph.insertParagraph(line, \"Before\");\n}\n\nawait context.sync();\n", "ps1": "# Add the new \"to-do\" lines after \"Ajouter doc avec licence et nos noms\",\n# keeping the existing trailing empty paragraph at the very end of the body.\n$d = $word.ActiveDocument\n\n# Helper: XML-escape plain text for a <w:t> element.\nfunction Escape-WordXml([string]$s) {\n    $s.Replace('&', '&amp;').Replace('<', '&lt;').Replace('>', '&gt;')\n}\n\n$wNs = 'xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\"'\n\n$newLines = @(\n    \"Mettre une scrollbarre sur la page principale de l\u2019app\",\n    \"Modifier la livraison quand on modifie le produit\",\n    \"\",\n    \"Questions : \",\n    \"\",\n    \"Creer objet prime\",\n    \"Creer objet arret maladie\",\n    \"Page avec tous les employes et on peut completer avec les salaires qui ont ete modifies\",\n    \"Charges ou les impactees\",\n    \"\"\n)\n\n$fragment = \"\"\nforeach ($line in $newLines) {\n    if ([string]::IsNullOrEmpty($line)) {\n        $fragment += \"<w:p $wNs/>\"\n    } else {\n        $escaped = Escape-WordXml $line\n        $fragment += \"<w:p $wNs><w:r><w:t xml:space=`\"preserve`\">$escaped</w:t></w:r></w:p>\"\n    }\n}\n\n# The document currently ends with one empty paragraph. Replace that\n# paragraph's whole range with: new paragraphs + a fresh trailing empty\n# paragraph, so the final empty paragraph is preserved at the end.\n$count = $d.Paragraphs.Count\n$lastParagraph = $d.Paragraphs.Item($count)\n$targetRange = $lastParagraph.Range\n$targetRange.InsertXML($fragment)\n"}
